$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.863.04"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.441.55"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "574.51"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "159.23"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.441.00"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  -5.49%  "
$ws.Range("D10").Value = "7.21"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").Value = "0.122"
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "4.034.73"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "27.73"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "64.886.76"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "3.405.47"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "382.21"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").Value = "7.97"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "72.21"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").Value = "9.84"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -4.01%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "161.12"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "2.910.91"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "0.0748"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").Value = "6.71"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("D41").Value = "26.28"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("D43").Value = "42.96"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "26.01"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.26"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "317.31"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  -4.62%  "
